$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where only F (Protect the Pine Trees - Code score) gets a perfect 10
$fullRows = @(8,9,10,11,12,14,15,16,17,21,22)
foreach ($r in $fullRows) {
    $ws.Range("F$r").Value = 10
}

# Row 22 also gets the Testcase score (G) filled in
$ws.Range("G22").Value = 10

# Rows with a partial score (9.75) and a remark ("no comment")
$partialRows = @(13,18,20)
foreach ($r in $partialRows) {
    $ws.Range("E$r").Value = "no comment"
    $ws.Range("F$r").Value = 9.75
}

# Reflect the reviewer's final cursor position / zoom level
$ws.Activate()
$ws.Range("E25").Select() | Out-Null
$excel.ActiveWindow.Zoom = 93
